$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old last row (372), pushing it down to 374.
$ws.Rows("372:373").Insert()

# Copy the date-cell format (column A) down into the two new rows so they
# keep the same style as the rest of the data (s="2").
$ws.Range("A371").Copy()
$ws.Range("A372:A373").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Ensure symbol column for the two new rows carries the same label.
$ws.Range("B372").Value2 = "ECONOMICS:PEM2"
$ws.Range("B373").Value2 = "ECONOMICS:PEM2"

# Row 368: value update only (date/symbol unchanged).
$ws.Range("C368").Value2 = 192778000000
$ws.Range("D368").Value2 = 192778000000
$ws.Range("E368").Value2 = 192778000000
$ws.Range("F368").Value2 = 192778000000

# Row 369: new date + new values.
$ws.Range("A369").Value2 = 44774.41666666666
$ws.Range("C369").Value2 = 295831000000
$ws.Range("D369").Value2 = 295831000000
$ws.Range("E369").Value2 = 295831000000
$ws.Range("F369").Value2 = 295831000000

# Row 370: new date + new values.
$ws.Range("A370").Value2 = 44805.41666666666
$ws.Range("C370").Value2 = 295831000000
$ws.Range("D370").Value2 = 295831000000
$ws.Range("E370").Value2 = 295831000000
$ws.Range("F370").Value2 = 295831000000

# Row 371: new date, values become the old row-369 values.
$ws.Range("A371").Value2 = 44835.41666666666
$ws.Range("C371").Value2 = 292818000000
$ws.Range("D371").Value2 = 292818000000
$ws.Range("E371").Value2 = 292818000000
$ws.Range("F371").Value2 = 292818000000

# Row 372 (newly inserted): old row-370 data.
$ws.Range("A372").Value2 = 44866.45833333334
$ws.Range("C372").Value2 = 288445000000
$ws.Range("D372").Value2 = 288445000000
$ws.Range("E372").Value2 = 288445000000
$ws.Range("F372").Value2 = 288445000000
$ws.Range("G372").Value2 = 0

# Row 373 (newly inserted): old row-371 data.
$ws.Range("A373").Value2 = 44896.45833333334
$ws.Range("C373").Value2 = 196800000000
$ws.Range("D373").Value2 = 196800000000
$ws.Range("E373").Value2 = 196800000000
$ws.Range("F373").Value2 = 196800000000
$ws.Range("G373").Value2 = 0

# Row 374 keeps the old row-372 values/date (untouched by the insert's shift).
